$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Define the name "pole" pointing to data!$E$30
$wb.Names.Add("pole", "=data!`$E`$30")

# Update the F column formulas to reference the defined name instead of E$30
$ws.Range("F2").Formula = '=IF(pole>E2,"ne","více")'
$ws.Range("F3:F23").Formula = '=IF(pole>E3,"ne","více")'

# Clear the J3:K5 values and remove the L3:L5 array formula / values
$ws.Range("J3:K5").Value = $null
$ws.Range("L3:L5").Value = $null

# Update selection to I11 as shown in the diff
$ws.Range("I11").Select()
